$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.114.47"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").Value = "1.750.74"
$ws.Range("E3").Value = "  +0.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.81"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5284"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2810"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06199"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.67%  "

$ws.Range("D10").Value = "1.746.11"
$ws.Range("E10").Value = "  +0.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07179"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.52"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6477"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.628"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "78.56"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9993"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("D18").Value = "26.014.27"
$ws.Range("E18").Value = "  +0.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006746"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.42%  "

$ws.Range("D21").Value = "1.969.48"
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.339"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.756"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.254"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.66%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.97"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.519"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("E27").Value = "  +2.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.817"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "104.82"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08295"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.816"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.664"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04604"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.10%  "

$ws.Range("E34").Value = "  +1.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.016"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6371"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.709"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01609"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.984"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.04%  "

$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.69"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3959"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7461"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.033"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1154"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.396"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05344"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.41"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.97"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3481"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.584"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.31%  "
